$d = $word.ActiveDocument

# --- 1. Remove the hidden "_GoBack" bookmark left over from the previous edit session ---
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
    # bookmark not present / not accessible - nothing to do
}

# --- 2. Insert a blank paragraph between the "Output" explanation paragraph and the table ---
$outputPara = $d.Paragraphs.Item(9)
$outputPara.Range.InsertParagraphAfter()

# --- 3. Append the new AND/OR warning paragraphs after the table ---
$lines = @(
    "AND/OR",
    "",
    "Varsel: Over 90% av journalenhetene er av samme type/navn.",
    "",
    "AND/OR",
    "",
    "Ingen journalenheter funnet."
)

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertRange = $lastPara.Range

foreach ($line in $lines) {
    $insertRange.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    if ($line -ne "") {
        $newPara.Range.Text = $line
    }
    $insertRange = $newPara.Range
}

Write-Output ("Final Paragraphs.Count=" + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output ("[$i] '" + $p.Range.Text + "'")
}
